# Auto-generated edit script: appends rows 52-57 to the "Artfynd" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52
$ws.Range("A52").Value = 112547017
$ws.Range("B52").Value = 90800
$ws.Range("C52").Value = "Ovaliderad"
$ws.Range("D52").Value = "NT"
$ws.Range("E52").Value = 3100
$ws.Range("F52").Value = "Talltaggsvamp"
$ws.Range("G52").Value = "Bankera fuligineoalba"
$ws.Range("H52").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("P52").Value = "Gustavbacke (Gustavbacke), Jmt"
$ws.Range("Q52").Value = 439919
$ws.Range("R52").Value = 6952481
$ws.Range("S52").Value = 25
$ws.Range("T52").Value = "Jämtland"
$ws.Range("U52").Value = "Härjedalen"
$ws.Range("V52").Value = "Jämtland"
$ws.Range("W52").Value = "Vemdalen"
$ws.Range("Y52").Value = "'2023-10-06"
$ws.Range("Y52").Style = "Normal"
$ws.Range("Z52").Value = "08:03"
$ws.Range("AA52").Value = "'2023-10-06"
$ws.Range("AA52").Style = "Normal"
$ws.Range("AB52").Value = "08:03"
$ws.Range("AD52").Value = $false
$ws.Range("AE52").Value = $false
$ws.Range("AG52").Value = $false
$ws.Range("AW52").Value = "Håkan Blomqvist"
$ws.Range("AX52").Value = "Håkan Blomqvist"

# Row 53
$ws.Range("A53").Value = 112546927
$ws.Range("B53").Value = 90800
$ws.Range("C53").Value = "Ovaliderad"
$ws.Range("D53").Value = "NT"
$ws.Range("E53").Value = 3100
$ws.Range("F53").Value = "Talltaggsvamp"
$ws.Range("G53").Value = "Bankera fuligineoalba"
$ws.Range("H53").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("P53").Value = "Gustavbacke, Jmt"
$ws.Range("Q53").Value = 440071
$ws.Range("R53").Value = 6952367
$ws.Range("S53").Value = 25
$ws.Range("T53").Value = "Jämtland"
$ws.Range("U53").Value = "Härjedalen"
$ws.Range("V53").Value = "Jämtland"
$ws.Range("W53").Value = "Vemdalen"
$ws.Range("Y53").Value = "'2023-10-06"
$ws.Range("Y53").Style = "Normal"
$ws.Range("Z53").Value = "07:50"
$ws.Range("AA53").Value = "'2023-10-06"
$ws.Range("AA53").Style = "Normal"
$ws.Range("AB53").Value = "07:50"
$ws.Range("AD53").Value = $false
$ws.Range("AE53").Value = $false
$ws.Range("AG53").Value = $false
$ws.Range("AW53").Value = "Håkan Blomqvist"
$ws.Range("AX53").Value = "Håkan Blomqvist"

# Row 54
$ws.Range("A54").Value = 112547097
$ws.Range("B54").Value = 90800
$ws.Range("C54").Value = "Ovaliderad"
$ws.Range("D54").Value = "NT"
$ws.Range("E54").Value = 3100
$ws.Range("F54").Value = "Talltaggsvamp"
$ws.Range("G54").Value = "Bankera fuligineoalba"
$ws.Range("H54").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("P54").Value = "Gustavbacke (Gustavbacke), Jmt"
$ws.Range("Q54").Value = 439879
$ws.Range("R54").Value = 6952505
$ws.Range("S54").Value = 25
$ws.Range("T54").Value = "Jämtland"
$ws.Range("U54").Value = "Härjedalen"
$ws.Range("V54").Value = "Jämtland"
$ws.Range("W54").Value = "Vemdalen"
$ws.Range("Y54").Value = "'2023-10-06"
$ws.Range("Y54").Style = "Normal"
$ws.Range("Z54").Value = "08:11"
$ws.Range("AA54").Value = "'2023-10-06"
$ws.Range("AA54").Style = "Normal"
$ws.Range("AB54").Value = "08:11"
$ws.Range("AD54").Value = $false
$ws.Range("AE54").Value = $false
$ws.Range("AG54").Value = $false
$ws.Range("AW54").Value = "Håkan Blomqvist"
$ws.Range("AX54").Value = "Håkan Blomqvist"

# Row 55
$ws.Range("A55").Value = 112547159
$ws.Range("B55").Value = 90830
$ws.Range("C55").Value = "Ovaliderad"
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 2059
$ws.Range("F55").Value = "Skrovlig taggsvamp"
$ws.Range("G55").Value = "Hydnellum scabrosum"
$ws.Range("H55").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("P55").Value = "Gustavbacke (Gustavbacke), Jmt"
$ws.Range("Q55").Value = 439800
$ws.Range("R55").Value = 6952502
$ws.Range("S55").Value = 25
$ws.Range("T55").Value = "Jämtland"
$ws.Range("U55").Value = "Härjedalen"
$ws.Range("V55").Value = "Jämtland"
$ws.Range("W55").Value = "Vemdalen"
$ws.Range("Y55").Value = "'2023-10-06"
$ws.Range("Y55").Style = "Normal"
$ws.Range("Z55").Value = "08:14"
$ws.Range("AA55").Value = "'2023-10-06"
$ws.Range("AA55").Style = "Normal"
$ws.Range("AB55").Value = "08:14"
$ws.Range("AD55").Value = $false
$ws.Range("AE55").Value = $false
$ws.Range("AG55").Value = $false
$ws.Range("AW55").Value = "Håkan Blomqvist"
$ws.Range("AX55").Value = "Håkan Blomqvist"

# Row 56
$ws.Range("A56").Value = 112546970
$ws.Range("B56").Value = 90808
$ws.Range("C56").Value = "Ovaliderad"
$ws.Range("D56").Value = "NT"
$ws.Range("E56").Value = 4362
$ws.Range("F56").Value = "Blå taggsvamp"
$ws.Range("G56").Value = "Hydnellum caeruleum"
$ws.Range("H56").Value = "(Hornem.) P.Karst."
$ws.Range("P56").Value = "Gustavbacke (Gustavbacke), Jmt"
$ws.Range("Q56").Value = 439971
$ws.Range("R56").Value = 6952512
$ws.Range("S56").Value = 25
$ws.Range("T56").Value = "Jämtland"
$ws.Range("U56").Value = "Härjedalen"
$ws.Range("V56").Value = "Jämtland"
$ws.Range("W56").Value = "Vemdalen"
$ws.Range("Y56").Value = "'2023-10-06"
$ws.Range("Y56").Style = "Normal"
$ws.Range("Z56").Value = "08:01"
$ws.Range("AA56").Value = "'2023-10-06"
$ws.Range("AA56").Style = "Normal"
$ws.Range("AB56").Value = "08:01"
$ws.Range("AD56").Value = $false
$ws.Range("AE56").Value = $false
$ws.Range("AG56").Value = $false
$ws.Range("AW56").Value = "Håkan Blomqvist"
$ws.Range("AX56").Value = "Håkan Blomqvist"

# Row 57
$ws.Range("A57").Value = 112546997
$ws.Range("B57").Value = 90808
$ws.Range("C57").Value = "Ovaliderad"
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 4362
$ws.Range("F57").Value = "Blå taggsvamp"
$ws.Range("G57").Value = "Hydnellum caeruleum"
$ws.Range("H57").Value = "(Hornem.) P.Karst."
$ws.Range("P57").Value = "Gustavbacke (Gustavbacke), Jmt"
$ws.Range("Q57").Value = 439940
$ws.Range("R57").Value = 6952516
$ws.Range("S57").Value = 25
$ws.Range("T57").Value = "Jämtland"
$ws.Range("U57").Value = "Härjedalen"
$ws.Range("V57").Value = "Jämtland"
$ws.Range("W57").Value = "Vemdalen"
$ws.Range("Y57").Value = "'2023-10-06"
$ws.Range("Y57").Style = "Normal"
$ws.Range("Z57").Value = "08:03"
$ws.Range("AA57").Value = "'2023-10-06"
$ws.Range("AA57").Style = "Normal"
$ws.Range("AB57").Value = "08:03"
$ws.Range("AD57").Value = $false
$ws.Range("AE57").Value = $false
$ws.Range("AG57").Value = $false
$ws.Range("AW57").Value = "Håkan Blomqvist"
$ws.Range("AX57").Value = "Håkan Blomqvist"
